$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 4,17
$arr[0,0] = 15
$arr[0,1] = 1
$arr[0,2] = 33.124
$arr[0,3] = 495.824
$arr[0,4] = 1.76684556
$arr[0,5] = 31235921.16
$arr[0,6] = 0.09416868
$arr[0,7] = 0.20787342
$arr[0,8] = 7.18301338
$arr[0,9] = 6.918986106177331
$arr[0,10] = 103.5211779907144
$arr[0,11] = 0.3522973706837276
$arr[0,12] = 14600205.03556428
$arr[0,13] = 0.01862986839057378
$arr[0,14] = 0.05919403656478638
$arr[0,15] = 3.445140540961417
$arr[0,16] = 0.15
$arr[1,0] = 15
$arr[1,1] = 1
$arr[1,2] = 45.036
$arr[1,3] = 656.188
$arr[1,4] = 1.3126927
$arr[1,5] = 30581864.56
$arr[1,6] = 0.06812261999999999
$arr[1,7] = 0.204237
$arr[1,8] = 9.4420936
$arr[1,9] = 7.535533992033358
$arr[1,10] = 106.0260898054761
$arr[1,11] = 0.1974610677928634
$arr[1,12] = 10275691.9165597
$arr[1,13] = 0.00852820517474839
$arr[1,14] = 0.04144504982395147
$arr[1,15] = 3.546223959162438
$arr[1,16] = 0.85
$arr[2,0] = 20
$arr[2,1] = 1
$arr[2,2] = 24.656
$arr[2,3] = 491.87
$arr[2,4] = 1.78192656
$arr[2,5] = 49877230.024
$arr[2,6] = 0.13100712
$arr[2,7] = 0.16222742
$arr[2,8] = 4.189310320000001
$arr[2,9] = 5.041080140005306
$arr[2,10] = 100.5283492254597
$arr[2,11] = 0.3670499658998767
$arr[2,12] = 20348057.83743384
$arr[2,13] = 0.02335894581535754
$arr[2,14] = 0.04862271945127375
$arr[2,15] = 2.110765459261526
$arr[2,16] = 0.15
$arr[3,0] = 20
$arr[3,1] = 1
$arr[3,2] = 34.578
$arr[3,3] = 666.562
$arr[3,4] = 1.29426156
$arr[3,5] = 54204474.64
$arr[3,6] = 0.09939218000000001
$arr[3,7] = 0.17193766
$arr[3,8] = 6.11989512
$arr[3,9] = 5.826475059866513
$arr[3,10] = 108.3662893215535
$arr[3,11] = 0.2056288652139693
$arr[3,12] = 16340568.5136297
$arr[3,13] = 0.01072079960338182
$arr[3,14] = 0.0355325839782691
$arr[3,15] = 2.340855516793035
$arr[3,16] = 0.85

$ws.Range("A14:Q17").Value = $arr
